$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "3.483.31"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.44%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "67.307.06"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  +7.83%  "
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "2.832.82"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -2.49%  "
